$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New dummy transaction row (row 7): HUDCO buy
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 45601
$ws.Cells.Item(7, 3).Value = "HUDCO"
$ws.Cells.Item(7, 4).Value = "HUDCO.NS"
$ws.Cells.Item(7, 5).Value = "Buy"
$ws.Cells.Item(7, 6).Value = 50
$ws.Cells.Item(7, 7).Value = 214.5
$ws.Cells.Item(7, 8).Formula = "=PRODUCT(F7,G7)"

# Move the active selection to the new total-amount cell
[void]$ws.Range("H7").Select()
